# Documentation Checklist update
#
# Marks "Alvin" as the documenter for four more source files
# (core/Stage.js [row 15], system/Canvas.js [row 38], sound/Sound.js
# [row 69] and sound/SoundManager.js [row 70]), and carries over the
# minor view/formatting state (scroll position, active selection, row
# heights, column widths) from the author's saved workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content: add "Alvin" in column B ("Documented by") ---
$ws.Range("B15").Value = "Alvin"
$ws.Range("B38").Value = "Alvin"
$ws.Range("B69").Value = "Alvin"
$ws.Range("B70").Value = "Alvin"

# --- Row heights: minor re-measurement of a few rows ---
$ws.Rows.Item(7).RowHeight = 14.75
$ws.Rows.Item(17).RowHeight = 14.75
$ws.Rows.Item(33).RowHeight = 14.75
$ws.Rows.Item(69).RowHeight = 14.9
$ws.Rows.Item(70).RowHeight = 14.9

# --- Column widths: small width bump across columns A:E ---
$ws.Columns.Item(1).ColumnWidth = 37.8313725490196
$ws.Columns.Item(2).ColumnWidth = 33.1294117647059
$ws.Columns.Item(3).ColumnWidth = 28.2627450980392
$ws.Columns.Item(4).ColumnWidth = 18.5607843137255
$ws.Columns.Item(5).ColumnWidth = 8.93725490196078

# --- View state: scroll position + active selection moved to B70 ---
$ws.Range("A51").Select() | Out-Null
$ws.Range("B70").Select() | Out-Null
